$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column U with header "Decommissioning date", matching the
# style of the existing last header cell (T1), and an empty U2 cell with
# the same style as T2.
$ws.Range("T1:T2").Copy()
$ws.Range("U1:U2").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("U1").Value = "Decommissioning date"

# Copy column width from column T (R:T share the same 23.5 width) to U.
$ws.Range("U1").ColumnWidth = $ws.Range("T1").ColumnWidth

# Remove the now-obsolete empty template rows (3-10) that only carried
# blank formatted cells.
$ws.Range("A3:T10").EntireRow.Delete()
